# "cambios antes de liberar" - Presupuesto Mundial / Torneo Chico updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Torneo Chico"
$ws.Activate()

# --- Row 3: Medalla de oro -> cantidad 1, total ahora es un valor fijo (50000) ---
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 50000

# --- Insert 4 new rows (6-9) for "Steam" prizes, right above the "Colación" block ---
$ws.Rows.Item(6).Resize(4).Insert()

$ws.Range("B6").Value = "Primer lugar"
$ws.Range("C6").Value = "Steam"
$ws.Range("D6").Value = 20000
$ws.Range("E6").Value = 1
$ws.Range("F6").Formula = "=D6*E6"

$ws.Range("B7").Value = "Segundo Lugar"
$ws.Range("C7").Value = "Steam"
$ws.Range("D7").Value = 60000
$ws.Range("E7").Value = 1
$ws.Range("F7").Formula = "=D7*E7"

$ws.Range("B8").Value = "Tercer Lugar"
$ws.Range("C8").Value = "Steam"
$ws.Range("D8").Value = 30000
$ws.Range("E8").Value = 1
$ws.Range("F8").Formula = "=D8*E8"

$ws.Range("B9").Value = "Copa de Plata"
$ws.Range("C9").Value = "Steam"
$ws.Range("D9").Value = 30000
$ws.Range("E9").Value = 1
$ws.Range("F9").Formula = "=D9*E9"

# --- Almuerzo block (now rows 16-22 after the insert above): quantities go to 0, ---
# --- "Bebidas" (row 22) no longer scales with the number of teams.              ---
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("E22").Value = 36

# --- Insert 3 new rows (23-25): "Pizzas" line item plus two blank spacer rows ---
$ws.Rows.Item(23).Resize(3).Insert()

$ws.Range("B23").Value = "Pizzas"
$ws.Range("C23").Value = "LaTorre"
$ws.Range("D23").Value = 5500
$ws.Range("E23").Value = 40
$ws.Range("F23").Formula = "=D23*E23"

# --- "Precio Final" per team goes up from 21000 to 25000 ---
$ws.Range("B38").Value = 25000

# --- Restore the on-screen selection/scroll position ---
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C16").Select()
